# Importer: prevent product to be linked to variation child
# Insert a new test row (row 46) into the product sample import sheet that
# exercises "parent sku points at a variation child" validation, pushing the
# previously-existing rows 46-50 down to 47-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 46 - shifts old rows 46:50 -> 47:51
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new test case data:
# SKU 42, Parent SKU 28, Name "This tries to link variation to child",
# variation values Color/Black + Size/XS, QTY 12, Price 100, same
# category/image columns as its neighbouring rows. No supplier on this row.
$ws.Cells.Item(46, 2).Value = 42
$ws.Cells.Item(46, 3).Value = 28
$ws.Cells.Item(46, 4).Value = "This tries to link variation to child"
$ws.Cells.Item(46, 6).Value = "Color/Black"
$ws.Cells.Item(46, 7).Value = "Size/XS"
$ws.Cells.Item(46, 9).Value = 12
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = "Test Category"
$ws.Cells.Item(46, 12).Value = "Test Category"
$ws.Cells.Item(46, 13).Value = "shirt1.jpeg"
$ws.Cells.Item(46, 14).Value = "shirt2.jpeg,shirt3.jpeg"

# Match the author's final selection state (cell B52 on the now-51-row sheet)
$ws.Range("B52").Select()
